# Rename the ontology term "budget area" -> "priority area" (and its
# Class/budgetArea/budgetSubArea variants) throughout the workbook.

$wb = $excel.ActiveWorkbook

# --- "Properties" sheet: budgetSubArea/budgetArea become prioritySubArea/priorityArea ---
$wsProps = $wb.Worksheets.Item("Properties")
$wsProps.Cells.Item(28, 1).Value = "prioritySubArea"
$wsProps.Cells.Item(28, 2).Value = "priority area subdivision"
$wsProps.Cells.Item(25, 1).Value = "priorityArea"
$wsProps.Cells.Item(25, 2).Value = "priority area"

# --- "Classes" sheet: the Area class becomes the PriorityArea class ---
$wsClasses = $wb.Worksheets.Item("Classes")
$wsClasses.Cells.Item(7, 1).Value = "PriorityArea"
$wsClasses.Cells.Item(7, 2).Value = "Priority Area"

# --- Restore the selections left behind on each sheet, and make "Classes" the active tab ---
$null = $wsProps.Range("B26").Select()
$null = $wsClasses.Activate()
$null = $wsClasses.Range("A8").Select()
